# Applies the changes described by the diff:
#  1. Insert a new row (004431546, GABRIELA, 60567.21) immediately above
#     the row for account 004254210 (MARCO).
#  2. Delete the row for account 005190138 (ANA, 7309).
#  3. Delete the row for account 004574428 (GUILHERME, 1524.94).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new GABRIELA row just above the MARCO row ---------------
$marcoRow = $ws.Columns.Item(1).Find("004254210").Row
$ws.Rows.Item($marcoRow).Insert()

$newRow = $marcoRow
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "004431546"
$ws.Cells.Item($newRow, 2).Value = "GABRIELA"
$ws.Cells.Item($newRow, 3).Value = 60567.21

# --- 2. Delete the ANA row ---------------------------------------------------
$anaRow = $ws.Columns.Item(1).Find("005190138").Row
$ws.Rows.Item($anaRow).Delete()

# --- 3. Delete the GUILHERME (004574428) row --------------------------------
$guilhermeRow = $ws.Columns.Item(1).Find("004574428").Row
$ws.Rows.Item($guilhermeRow).Delete()
